# Apply the "cryptos list" update (GitHub Actions scheduled refresh).
# Each row's Price (D) and/or Volume(1h) (E) cell is rewritten with the
# latest scraped value; row 51 swapped out one listing entirely.
#
# The Price column stores values as plain text (prices use "." as both a
# thousands separator and decimal point, e.g. "27.940.86", so Excel can't
# treat the column as numeric). Some new prices (e.g. "213.45") happen to
# look like valid numbers, and a plain .Value assignment would silently
# coerce them to floating point (losing the "text" cell type and picking
# up binary-float rounding noise). To keep those cells text - matching
# every other cell in the column - briefly force a text NumberFormat
# before writing the value, then restore the default "Normal" style so no
# stray per-cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($ref, $value) {
    $ws.Range($ref).Value = $value
}

function Set-TextCell($ref, $value) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $value
    $ws.Range($ref).Style = "Normal"
}

Set-Cell "D2"  "27.940.86"
Set-Cell "E2"  "  +1.48%  "

Set-Cell "D3"  "1.642.71"

Set-TextCell "D5" "213.45"
Set-Cell "E5"  "  +0.80%  "

Set-Cell "E6"  "  -0.16%  "

Set-Cell "E7"  "  -0.06%  "

Set-TextCell "D8" "23.94"
Set-Cell "E8"  "  +3.23%  "

Set-Cell "E9"  "  +0.65%  "

Set-TextCell "D11" "0.0877"
Set-Cell "E11" "  -0.39%  "

Set-Cell "D12" "1.876.07"
Set-Cell "E12" "  +1.19%  "

Set-Cell "D13" "1.643.12"
Set-Cell "E13" "  +0.81%  "

Set-TextCell "D14" "0.577"
Set-Cell "E14" "  +4.91%  "

Set-Cell "E15" "  +1.03%  "

Set-TextCell "D16" "65.94"
Set-Cell "E16" "  +1.20%  "

Set-Cell "D17" "27.916.54"
Set-Cell "E17" "  +1.52%  "

Set-TextCell "D18" "230.96"
Set-Cell "E18" "  +0.36%  "

Set-Cell "E19" "  +0.97%  "

Set-TextCell "D20" "7.64"
Set-Cell "E20" "  +1.42%  "

Set-Cell "E21" "  -0.04%  "

Set-Cell "E22" "  +7.08%  "

Set-Cell "E23" "  +1.58%  "

Set-Cell "E24" "  -0.21%  "

Set-TextCell "D25" "152.43"
Set-Cell "E25" "  +2.42%  "

Set-TextCell "D26" "6.95"
Set-Cell "E26" "  +0.91%  "

Set-Cell "E27" "  +0.79%  "

Set-TextCell "D28" "15.75"
Set-Cell "E28" "  +1.36%  "

Set-Cell "E29" "  +0.00%  "

Set-Cell "E30" "  +1.14%  "

Set-Cell "E32" "  +2.19%  "

Set-Cell "D33" "1.423.61"
Set-Cell "E33" "  -2.99%  "

Set-Cell "E34" "  +2.22%  "

Set-Cell "E35" "  +2.19%  "

Set-Cell "E36" "  +0.18%  "

Set-Cell "E37" "  +2.31%  "

Set-TextCell "D38" "0.927"
Set-Cell "E38" "  -0.47%  "

Set-Cell "E39" "  +1.02%  "

Set-Cell "E40" "  +0.68%  "

Set-Cell "E41" "  +2.04%  "

Set-Cell "E43" "  +0.18%  "

Set-Cell "E44" "  +0.29%  "

Set-TextCell "D45" "5.46"
Set-Cell "E45" "  +3.25%  "

Set-Cell "E46" "  +3.59%  "

Set-Cell "E47" "  +0.23%  "

Set-Cell "D48" "1.784.37"

Set-Cell "E49" "  +1.85%  "

Set-Cell "E50" "  +0.96%  "

# Row 51: EnergySwap -> Cronos (listing swapped to a different coin)
Set-Cell "B51" "Cronos"
Set-Cell "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D51" "0.0506"
Set-Cell "E51" "  +0.52%  "
